$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 73.58335264412017
$ws.Range("B3").Value = 0.916654408844484
$ws.Range("B4").Value = 0.05921143599821604
$ws.Range("B5").Value = 0.3707876469228513
